# Update the Gaziantep Buyuksehir Belediyesi website URL (row 12, column B)
# from the old "www.gaziantep.bel.tr" to the new scraped value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "https://www.gaziantep.bel.tr/tr"

# Widen column B (Web) slightly to fit the longer URL.
# ColumnWidth is in "character" units; Excel's saved <col width="..">
# attribute adds a fixed ~0.83 padding offset on top of this value, so we
# back that off here to land exactly on width="28" in the saved XML
# (previously width="26").
$ws.Columns.Item(2).ColumnWidth = 27.166666666666668

# The author's last action left the selection on B2.
$ws.Range("B2").Select()
